$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 5: Node.JS -----------------------------------------------------
# Ring moves from Adopt to Hold
$ws.Range("B5").Value = "Hold"
# Description text is updated (quotes around the Node.js® sentence removed)
$ws.Range("E5").Value = "Node.js® is a JavaScript runtime built on Chrome's V8 JavaScript engine. JS on the backend. "
# Vendor is now filled in
$ws.Range("F5").Value = "Node"

# --- Row 6: Rust ----------------------------------------------------------
# Ring moves from Assess to Adopt
$ws.Range("B6").Value = "Adopt"
# Description text updated (quotes removed, trailing period removed)
$ws.Range("E6").Value = "Rust is blazingly fast and memory-efficient: with no runtime or garbage collector, it can power performance-critical services, run on embedded devices, and easily integrate with other languages"
# Vendor stays Mozilla
$ws.Range("F6").Value = "Mozilla"
# Assessed By is now filled in
$ws.Range("G6").Value = "Saša Slankamenac"
# Last Assessment Date is now filled in
$ws.Range("H6").Value = "Q3 - 2022"
# Customer Demand is now filled in
$ws.Range("I6").Value = "Limited"
# Should we adopt changes from No to Yes
$ws.Range("J6").Value = "Yes"
# Value is now filled in
$ws.Range("K6").Value = "mature enough to be used in micro-service architectures for web API micro service deployments as container orchestrated pods in k8s or serverless eg AWS lambda functions where processing performance and memory management are critical"

$wb.Save()
